$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

$assigned = @(
    'Jack Mogelof (9172163044) (D), Alejandro L (2676633214)',
    'Josh Greene (9176368454) (D), Thor Waguespack (2404299065) (D), Ali Awada (6469771844)',
    'Matheo (6462079196) (D)',
    'Jaxon (4159881691) (D), Jake Dieterich (9177968225), Blake Steel (6109995081) (D)',
    'Alexander (9175433831) (D), Henry (6073399363), Harry Corbin (5182223247) (D)',
    'Alejandro Espinosa (7866702380) (D), Jamari Pitchford (7738297627) (D), Adi (3103038889) (D)',
    'George Ryckman (9176015863) (D), Alejandro E. Ulvert (3058421968) (D), Noah Yaffe (6109995081) (D)',
    'Alexander (9175433831) (D), Henry (6073399363)',
    'Thor Waguespack (2404299065) (D), Edu (7862012521) (D), Adi (3103038889) (D)',
    'Jack Mogelof (9172163044) (D), Ben Kairouz (9179401653) (D)',
    'Ezana (2405594003), George Ryckman (9176015863) (D)',
    'Jake Dieterich (9177968225), Alejandro E. Ulvert (3058421968) (D), Noah Yaffe (6109995081) (D)',
    'Jaxon (4159881691) (D), Paul (3106583636) (D), Alejandro Espinosa (7866702380) (D)',
    'Jake Dieterich (9177968225), Adi (3103038889) (D)',
    'Paul (3106583636) (D), Josh Greene (9176368454) (D), Ali Awada (6469771844)',
    'Jaxon (4159881691) (D), Matheo (6462079196) (D), Blake Steel (6109995081) (D)',
    'Henry (6073399363), Edu (7862012521) (D)',
    'Ezana (2405594003), Jamari Pitchford (7738297627) (D), Harry Corbin (5182223247) (D)',
    'Kamsi (6465933923), Alejandro Espinosa (7866702380) (D), Noah Yaffe (6109995081) (D)',
    'Gabe Heller (7187220895) (D), Alexander (9175433831) (D), Alejandro L (2676633214)',
    'Ezana (2405594003), Matheo (6462079196) (D)',
    'Paul (3106583636) (D), Thor Waguespack (2404299065) (D), Edu (7862012521) (D)',
    'Jack Mogelof (9172163044) (D)',
    'Ben Kairouz (9179401653) (D), Kamsi (6465933923)',
    'Gabe Heller (7187220895) (D), Kamsi (6465933923)',
    'Jaxon (4159881691) (D), Paul (3106583636) (D), Alejandro L (2676633214)',
    'Jack Mogelof (9172163044) (D), Jamari Pitchford (7738297627) (D), Edu (7862012521) (D)',
    'Ben Kairouz (9179401653) (D), Thor Waguespack (2404299065) (D), Alejandro E. Ulvert (3058421968) (D)',
    'Ali Awada (6469771844), Blake Steel (6109995081) (D), Harry Corbin (5182223247) (D)',
    'Gabe Heller (7187220895) (D), Harry Corbin (5182223247) (D)',
    'Jaxon (4159881691) (D), Alejandro L (2676633214), Alejandro Espinosa (7866702380) (D)',
    'Josh Greene (9176368454) (D), Ali Awada (6469771844), Alejandro E. Ulvert (3058421968) (D)',
    'Ben Kairouz (9179401653) (D), Paul (3106583636) (D), George Ryckman (9176015863) (D)',
    'Jamari Pitchford (7738297627) (D), Jake Dieterich (9177968225), Blake Steel (6109995081) (D)',
    'Josh Greene (9176368454) (D), George Ryckman (9176015863) (D)',
    'Kamsi (6465933923), Alexander (9175433831) (D), Alejandro Espinosa (7866702380) (D)',
    'Matheo (6462079196) (D), Edu (7862012521) (D), Ali Awada (6469771844)',
    'Ben Kairouz (9179401653) (D), Henry (6073399363), Thor Waguespack (2404299065) (D)',
    'Jamari Pitchford (7738297627) (D), Blake Steel (6109995081) (D), Noah Yaffe (6109995081) (D)',
    'Ezana (2405594003), Alejandro E. Ulvert (3058421968) (D), Adi (3103038889) (D)'
)

for ($i = 0; $i -lt $assigned.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $assigned[$i]
}

